# Scheduled-runner refresh of cached currentAveragePrice / LevePrice /
# LeveProfit figures (columns H-N) across all eight Sheets tabs.
# Values come from a refreshed market-data snapshot; item/leve metadata
# (columns A-G) is untouched. Some rows gain/lose the LeveProfit* cells
# entirely when the corresponding Leve price is zero/nonzero.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1796.9166
$ws.Range("J17").Value = 1863.6364
$ws.Range("L17").Value = 5590.9092
$ws.Range("N17").Value = -5926.9092
$ws.Range("H80").Value = 521.55554
$ws.Range("I80").Value = 434.16666
$ws.Range("J80").Value = 696.3333
$ws.Range("K80").Value = 1302.49998
$ws.Range("L80").Value = 2088.9999
$ws.Range("M80").Value = -304.4999800000001
$ws.Range("N80").Value = -4084.9999
$ws.Range("H83").Value = 521.55554
$ws.Range("I83").Value = 434.16666
$ws.Range("J83").Value = 696.3333
$ws.Range("K83").Value = 3907.49994
$ws.Range("L83").Value = 6266.9997
$ws.Range("M83").Value = 1084.50006
$ws.Range("N83").Value = -16250.9997
$ws.Range("H100").Value = 2998.5
$ws.Range("I100").Value = 2331.3333
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 2331.3333
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -1790.3333
$ws.Range("N100").Value = -6082
$ws.Range("H137").Value = 1425.375
$ws.Range("I137").Value = 1242.1666
$ws.Range("K137").Value = 3726.4998
$ws.Range("M137").Value = -1176.4998
$ws.Range("H141").Value = 2481.1428
$ws.Range("I141").Value = 1518.1923
$ws.Range("K141").Value = 4554.5769
$ws.Range("M141").Value = 625.4231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2601.9092
$ws.Range("I45").Value = 2601.9092
$ws.Range("K45").Value = 2601.9092
$ws.Range("M45").Value = -2224.9092
$ws.Range("H61").Value = 2228.1538
$ws.Range("I61").Value = 2228.1538
$ws.Range("K61").Value = 2228.1538
$ws.Range("M61").Value = -2016.1538
$ws.Range("H74").Value = 904.8946999999999
$ws.Range("I74").Value = 871.8333
$ws.Range("K74").Value = 871.8333
$ws.Range("M74").Value = 2.166699999999992
$ws.Range("H77").Value = 904.8946999999999
$ws.Range("I77").Value = 871.8333
$ws.Range("K77").Value = 4359.1665
$ws.Range("M77").Value = 8.833499999999731
$ws.Range("H97").Value = 1076.7142
$ws.Range("I97").Value = 1056.1666
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 1056.1666
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -560.1666
$ws.Range("N97").Value = -2192
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = $null
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = 0
$ws.Range("H132").Value = 790.8461
$ws.Range("I132").Value = 790.8461
$ws.Range("K132").Value = 2372.5383
$ws.Range("M132").Value = 157.4616999999998
$ws.Range("H136").Value = 2228.1538
$ws.Range("I136").Value = 2228.1538
$ws.Range("K136").Value = 6684.4614
$ws.Range("M136").Value = -4134.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 37500
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -37122
$ws.Range("H84").Value = 37500
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -115608
$ws.Range("H99").Value = 2490
$ws.Range("I99").Value = 2157.7778
$ws.Range("J99").Value = 2822.2222
$ws.Range("K99").Value = 2157.7778
$ws.Range("L99").Value = 2822.2222
$ws.Range("M99").Value = -659.7777999999998
$ws.Range("N99").Value = -5818.2222
$ws.Range("H107").Value = 2979.2222
$ws.Range("I107").Value = 2974
$ws.Range("K107").Value = 2974
$ws.Range("M107").Value = -1054
$ws.Range("H109").Value = 72998.5
$ws.Range("J109").Value = 72998.5
$ws.Range("L109").Value = 72998.5
$ws.Range("N109").Value = -75772.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3380.9412
$ws.Range("I31").Value = 1791
$ws.Range("K31").Value = 1791
$ws.Range("M31").Value = -1496
$ws.Range("H34").Value = 3380.9412
$ws.Range("I34").Value = 1791
$ws.Range("K34").Value = 1791
$ws.Range("M34").Value = -1589
$ws.Range("H58").Value = 2852.0527
$ws.Range("I58").Value = 1065.5454
$ws.Range("K58").Value = 1065.5454
$ws.Range("M58").Value = -862.5454
$ws.Range("H132").Value = 4266
$ws.Range("I132").Value = 4449.5
$ws.Range("J132").Value = 3899
$ws.Range("K132").Value = 13348.5
$ws.Range("L132").Value = 11697
$ws.Range("M132").Value = -10818.5
$ws.Range("N132").Value = -16757
$ws.Range("H133").Value = 92985
$ws.Range("J133").Value = 92985
$ws.Range("L133").Value = 92985
$ws.Range("N133").Value = -98045
$ws.Range("H134").Value = 3233.9473
$ws.Range("I134").Value = 3222.7273
$ws.Range("K134").Value = 9668.1819
$ws.Range("M134").Value = -7133.1819
$ws.Range("H136").Value = 2852.0527
$ws.Range("I136").Value = 1065.5454
$ws.Range("K136").Value = 3196.6362
$ws.Range("M136").Value = -646.6361999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1328761.9
$ws.Range("J4").Value = 825
$ws.Range("L4").Value = 2475
$ws.Range("N4").Value = -2699
$ws.Range("H36").Value = 1699.875
$ws.Range("J36").Value = 799.5
$ws.Range("L36").Value = 2398.5
$ws.Range("N36").Value = -2736.5
$ws.Range("H133").Value = 4799.6665
$ws.Range("J133").Value = 11999
$ws.Range("L133").Value = 35997
$ws.Range("N133").Value = -46117
$ws.Range("H136").Value = 16331.667
$ws.Range("I136").Value = 9000
$ws.Range("J136").Value = 19997.5
$ws.Range("K136").Value = 27000
$ws.Range("L136").Value = 59992.5
$ws.Range("M136").Value = -21900
$ws.Range("N136").Value = -70192.5
$ws.Range("H139").Value = 7766
$ws.Range("J139").Value = 7766
$ws.Range("L139").Value = 23298
$ws.Range("N139").Value = -33578
$ws.Range("H141").Value = 9007.25
$ws.Range("J141").Value = 12000
$ws.Range("L141").Value = 36000
$ws.Range("N141").Value = -46360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3158.5625
$ws.Range("I102").Value = 2268.7
$ws.Range("K102").Value = 2268.7
$ws.Range("M102").Value = -646.6999999999998
$ws.Range("H116").Value = 70000
$ws.Range("J116").Value = 70000
$ws.Range("L116").Value = 70000
$ws.Range("N116").Value = -79178
$ws.Range("H122").Value = 146108.72
$ws.Range("I122").Value = 2698.25
$ws.Range("J122").Value = 337322.66
$ws.Range("K122").Value = 8094.75
$ws.Range("L122").Value = 1011967.98
$ws.Range("M122").Value = -5644.75
$ws.Range("N122").Value = -1016867.98
$ws.Range("H132").Value = 2076.125
$ws.Range("I132").Value = 1371.8
$ws.Range("K132").Value = 4115.4
$ws.Range("M132").Value = -1585.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2538.8333
$ws.Range("I16").Value = 2309.5
$ws.Range("J16").Value = 2997.5
$ws.Range("K16").Value = 2309.5
$ws.Range("L16").Value = 2997.5
$ws.Range("M16").Value = -2139.5
$ws.Range("N16").Value = -3337.5
$ws.Range("H39").Value = 46122.715
$ws.Range("I39").Value = 33929.5
$ws.Range("J39").Value = 51000
$ws.Range("K39").Value = 33929.5
$ws.Range("L39").Value = 51000
$ws.Range("M39").Value = -33469.5
$ws.Range("N39").Value = -51920
$ws.Range("H40").Value = 2609.111
$ws.Range("I40").Value = 2609.111
$ws.Range("K40").Value = 2609.111
$ws.Range("M40").Value = -2473.111
$ws.Range("H132").Value = 2793.2778
$ws.Range("I132").Value = 2442.1667
$ws.Range("K132").Value = 7326.500100000001
$ws.Range("M132").Value = -4796.500100000001
$ws.Range("H136").Value = 6246.7
$ws.Range("I136").Value = 5710.125
$ws.Range("K136").Value = 17130.375
$ws.Range("M136").Value = -14580.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 350000
$ws.Range("I26").Value = 10000
$ws.Range("J26").Value = 418000
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 418000
$ws.Range("M26").Value = -9707
$ws.Range("N26").Value = -418586
$ws.Range("H39").Value = 2999
$ws.Range("I39").Value = 2999
$ws.Range("K39").Value = 2999
$ws.Range("M39").Value = -2586
$ws.Range("H42").Value = 30049
$ws.Range("J42").Value = 30049
$ws.Range("L42").Value = 30049
$ws.Range("N42").Value = -30805
$ws.Range("H43").Value = 1500
$ws.Range("I43").Value = 1500
$ws.Range("K43").Value = 1500
$ws.Range("M43").Value = -1351
$ws.Range("H122").Value = 2382
$ws.Range("I122").Value = 2382
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7146
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -4696
$ws.Range("H132").Value = 3591.1
$ws.Range("I132").Value = 3217.1538
$ws.Range("J132").Value = 4285.5713
$ws.Range("K132").Value = 9651.4614
$ws.Range("L132").Value = 12856.7139
$ws.Range("M132").Value = -7121.4614
$ws.Range("N132").Value = -17916.7139
$ws.Range("H139").Value = 178914
$ws.Range("J139").Value = 178914
$ws.Range("L139").Value = 178914
$ws.Range("N139").Value = -189194
